# 16.6.1.xlsx update: add 2018/2019 "факт" columns, add 2020/2021 утв./факт columns,
# refresh deviation percentages and underlying data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend formatting into the new columns (AK:AM) by copying from the
#        last existing column (AJ) in the relevant rows, so the new cells
#        inherit the right style (border/font/number format). ---
$ws.Range("AJ2").Copy()
$ws.Range("AK2:AM2").PasteSpecial(-4122)

$ws.Range("AJ3").Copy()
$ws.Range("AK3:AM3").PasteSpecial(-4122)

$ws.Range("AJ5").Copy()
$ws.Range("AK5:AM11").PasteSpecial(-4122)

$ws.Range("AJ12").Copy()
$ws.Range("AK12:AM12").PasteSpecial(-4122)

# --- 2. Update the header row (row 3) text. ---
$ws.Range("AC3").Value = "2018 факт "
$ws.Range("AE3").Value = "2019 утв."
$ws.Range("AF3").Value = "2019 факт"
$ws.Range("AH3").Value = "2020 утв."
$ws.Range("AI3").Value = "2020 факт"
$ws.Range("AJ3").Value = "откл. от утв., %"
$ws.Range("AK3").Value = "2021 утв."
$ws.Range("AL3").Value = "2021 факт"
$ws.Range("AM3").Value = "откл. от утв., %"

# --- 3. Update the numeric data (rows 5-12, columns AE:AM). ---

# Row 5 - Жалпы мамлекеттик кызматтар / Государственные услуги общего назначения
$ws.Range("AE5").Value = 43737.8
$ws.Range("AF5").Value = 43258.3
$ws.Range("AG5").Value = 98.9
$ws.Range("AH5").Value = 46293.5
$ws.Range("AI5").Value = 47153.5
$ws.Range("AJ5").Value = 101.9
$ws.Range("AK5").Value = 47483.3
$ws.Range("AL5").Value = 52020.5
$ws.Range("AM5").Value = 109.6

# Row 6 - Экономикалык суроолор / Экономические вопросы
$ws.Range("AE6").Value = 6265.4
$ws.Range("AF6").Value = 4434.6
$ws.Range("AG6").Value = 70.8
$ws.Range("AH6").Value = 7935.8
$ws.Range("AI6").Value = 3895.8
$ws.Range("AJ6").Value = 49.1
$ws.Range("AK6").Value = 8997
$ws.Range("AL6").Value = 6212.4
$ws.Range("AM6").Value = 69

# Row 7 - Айлана чөйрөнү коргоо / Охрана окружающей среды
$ws.Range("AE7").Value = 728.5
$ws.Range("AF7").Value = 695.7
$ws.Range("AG7").Value = 95.5
$ws.Range("AH7").Value = 746.9
$ws.Range("AI7").Value = 583.2
$ws.Range("AJ7").Value = 78.1
$ws.Range("AK7").Value = 639.2
$ws.Range("AL7").Value = 600.8
$ws.Range("AM7").Value = 94

# Row 8 - Турак жай жана коммуналдык кызматтар / Жилищные и коммунальные услуги
$ws.Range("AE8").Value = 1249
$ws.Range("AF8").Value = 1244.7
$ws.Range("AG8").Value = 99.7
$ws.Range("AH8").Value = 1249
$ws.Range("AI8").Value = 1207.6
$ws.Range("AJ8").Value = 96.7
$ws.Range("AK8").Value = 1208.1
$ws.Range("AL8").Value = 1332.7
$ws.Range("AM8").Value = 110.3

# Row 9 - Саламаттыкты сактоо / Здравоохранение
$ws.Range("AE9").Value = 2582.6
$ws.Range("AF9").Value = 2477.5
$ws.Range("AG9").Value = 95.9
$ws.Range("AH9").Value = 3109
$ws.Range("AI9").Value = 3225.2
$ws.Range("AJ9").Value = 103.7
$ws.Range("AK9").Value = 3131.3
$ws.Range("AL9").Value = 4833.7
$ws.Range("AM9").Value = 154.4

# Row 10 - Эс алууну жана маданий-диний иш-чараларды уюштуруу / Организация отдыха и культурно-религиозная деятельность
$ws.Range("AE10").Value = 2686.4
$ws.Range("AF10").Value = 2829
$ws.Range("AG10").Value = 105.3
$ws.Range("AH10").Value = 2993.4
$ws.Range("AI10").Value = 2624.5
$ws.Range("AJ10").Value = 87.7
$ws.Range("AK10").Value = 2798.4
$ws.Range("AL10").Value = 3088
$ws.Range("AM10").Value = 110.3

# Row 11 - Билим берүү / Образование
$ws.Range("AE11").Value = 23397.4
$ws.Range("AF11").Value = 24364.8
$ws.Range("AG11").Value = 104.1
$ws.Range("AH11").Value = 30085.9
$ws.Range("AI11").Value = 29223.5
$ws.Range("AJ11").Value = 97.1
$ws.Range("AK11").Value = 30439.7
$ws.Range("AL11").Value = 30705.3
$ws.Range("AM11").Value = 100.9

# Row 12 - Соцалдык коргоо / Социальная защита
$ws.Range("AE12").Value = 13137.1
$ws.Range("AF12").Value = 10924.7
$ws.Range("AG12").Value = 83.2
$ws.Range("AH12").Value = 12158.7
$ws.Range("AI12").Value = 10980.3
$ws.Range("AJ12").Value = 90.3
$ws.Range("AK12").Value = 11664.9
$ws.Range("AL12").Value = 11939.1
$ws.Range("AM12").Value = 102.4

# --- 4. Update the selection (matches the saved view state in the workbook). ---
[void]$ws.Range("AF4").Select()

Write-Host "edit complete"
